$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data values in rows 2-5 (new sensor readings, timestamps shifted by +317 days)
$ws.Range("A2").Value = "45098.50694444445"
$ws.Range("B2").Value = "19.217"
$ws.Range("C2").Value = "12.901"
$ws.Range("D2").Value = "4.042"
$ws.Range("E2").Value = "40.812"
$ws.Range("F2").Value = "32.818"
$ws.Range("G2").Value = "15.123"
$ws.Range("H2").Value = "47.986"
$ws.Range("I2").Value = "23.269"
$ws.Range("J2").Value = "9.710000000000001"
$ws.Range("K2").Value = "14.67"
$ws.Range("L2").Value = "16.076"
$ws.Range("M2").Value = "16.742"
$ws.Range("N2").Value = "4.827"
$ws.Range("O2").Value = "15.038"
$ws.Range("P2").Value = "20.994"
$ws.Range("Q2").Value = "12.85"
$ws.Range("R2").Value = "3.46"
$ws.Range("S2").Value = "2.249"
$ws.Range("T2").Value = "221.547"
$ws.Range("U2").Value = "41.81"
$ws.Range("V2").Value = "13.881"
$ws.Range("W2").Value = "27.553"
$ws.Range("X2").Value = "14.055"
$ws.Range("Y2").Value = "3.03"
$ws.Range("Z2").Value = "24.312"
$ws.Range("AA2").Value = "12.261"
$ws.Range("AB2").Value = "11.125"
$ws.Range("AC2").Value = "13.047"
$ws.Range("AD2").Value = "16.565"
$ws.Range("AE2").Value = "3.456"
$ws.Range("AF2").Value = "42.557"
$ws.Range("AG2").Value = "7.647"
$ws.Range("AH2").Value = "17.354"

$ws.Range("A3").Value = "45098.51388888889"
$ws.Range("B3").Value = "1.441"
$ws.Range("C3").Value = "0.369"
$ws.Range("D3").Value = "1.308"
$ws.Range("E3").Value = "2.888"
$ws.Range("F3").Value = "1.794"
$ws.Range("G3").Value = "1.14"
$ws.Range("H3").Value = "11.081"
$ws.Range("I3").Value = "1.745"
$ws.Range("J3").Value = "0.612"
$ws.Range("K3").Value = "0.541"
$ws.Range("L3").Value = "1.118"
$ws.Range("M3").Value = "0.998"
$ws.Range("N3").Value = "0.399"
$ws.Range("O3").Value = "1.128"
$ws.Range("P3").Value = "1.568"
$ws.Range("Q3").Value = "1.346"
$ws.Range("R3").Value = "1.429"
$ws.Range("S3").Value = "0.623"
$ws.Range("T3").Value = "10.001"
$ws.Range("U3").Value = "3.735"
$ws.Range("V3").Value = "1.041"
$ws.Range("W3").Value = "2.258"
$ws.Range("X3").Value = "0.913"
$ws.Range("Y3").Value = "0.725"
$ws.Range("Z3").Value = "4.9"
$ws.Range("AA3").Value = "0.92"
$ws.Range("AB3").Value = "1.099"
$ws.Range("AC3").Value = "1.234"
$ws.Range("AD3").Value = "0.977"
$ws.Range("AE3").Value = "1.265"
$ws.Range("AF3").Value = "10.731"
$ws.Range("AG3").Value = "0.402"
$ws.Range("AH3").Value = "1.315"

$ws.Range("A4").Value = "45098.52083333334"
$ws.Range("B4").Value = "6.726"
$ws.Range("C4").Value = "4.589"
$ws.Range("D4").Value = "0.997"
$ws.Range("E4").Value = "14.515"
$ws.Range("F4").Value = "11.55"
$ws.Range("G4").Value = "5.293"
$ws.Range("H4").Value = "18.923"
$ws.Range("I4").Value = "8.144"
$ws.Range("J4").Value = "3.399"
$ws.Range("K4").Value = "4.989"
$ws.Range("L4").Value = "5.824"
$ws.Range("M4").Value = "6.03"
$ws.Range("N4").Value = "1.685"
$ws.Range("O4").Value = "5.263"
$ws.Range("P4").Value = "7.284"
$ws.Range("Q4").Value = "4.686"
$ws.Range("R4").Value = "0.969"
$ws.Range("S4").Value = "0.532"
$ws.Range("T4").Value = "72.736"
$ws.Range("U4").Value = "14.562"
$ws.Range("V4").Value = "4.858"
$ws.Range("W4").Value = "9.465999999999999"
$ws.Range("X4").Value = "4.97"
$ws.Range("Y4").Value = "1.09"
$ws.Range("Z4").Value = "8.989000000000001"
$ws.Range("AA4").Value = "4.291"
$ws.Range("AB4").Value = "3.962"
$ws.Range("AC4").Value = "4.633"
$ws.Range("AD4").Value = "5.999"
$ws.Range("AE4").Value = "0.773"
$ws.Range("AF4").Value = "16.844"
$ws.Range("AG4").Value = "2.622"
$ws.Range("AH4").Value = "6.073"

$ws.Range("A5").Value = "45098.52777777778"
$ws.Range("B5").Value = "6.73"
$ws.Range("C5").Value = "4.71"
$ws.Range("D5").Value = "0.78"
$ws.Range("E5").Value = "14.56"
$ws.Range("F5").Value = "11.68"
$ws.Range("G5").Value = "5.29"
$ws.Range("H5").Value = "20.29"
$ws.Range("I5").Value = "8.140000000000001"
$ws.Range("J5").Value = "3.49"
$ws.Range("K5").Value = "5.08"
$ws.Range("L5").Value = "5.85"
$ws.Range("M5").Value = "6.09"
$ws.Range("N5").Value = "1.69"
$ws.Range("O5").Value = "5.26"
$ws.Range("P5").Value = "7.38"
$ws.Range("Q5").Value = "4.62"
$ws.Range("R5").Value = "0.73"
$ws.Range("S5").Value = "0.43"
$ws.Range("T5").Value = "72.78"
$ws.Range("U5").Value = "14.72"
$ws.Range("V5").Value = "4.86"
$ws.Range("W5").Value = "9.69"
$ws.Range("X5").Value = "5.06"
$ws.Range("Y5").Value = "1"
$ws.Range("Z5").Value = "9.789999999999999"
$ws.Range("AA5").Value = "4.29"
$ws.Range("AB5").Value = "3.92"
$ws.Range("AC5").Value = "4.59"
$ws.Range("AD5").Value = "6.07"
$ws.Range("AE5").Value = "0.5600000000000001"
$ws.Range("AF5").Value = "18.36"
$ws.Range("AG5").Value = "2.65"
$ws.Range("AH5").Value = "6.07"

# Remove the last data row (row 6) - dataset now has one fewer sample
$ws.Rows("6:6").Delete()

# Adjust column widths (character units converted to Excel ColumnWidth offset of -0.83)
$ws.Columns("B").ColumnWidth = 7.17
$ws.Columns("C").ColumnWidth = 7.17
$ws.Columns("G").ColumnWidth = 7.17
$ws.Columns("I").ColumnWidth = 7.17
$ws.Columns("L").ColumnWidth = 7.17
$ws.Columns("M").ColumnWidth = 7.17
$ws.Columns("O").ColumnWidth = 7.17
$ws.Columns("P").ColumnWidth = 7.17
$ws.Columns("T").ColumnWidth = 8.17
$ws.Columns("V").ColumnWidth = 7.17
$ws.Columns("W").ColumnWidth = 7.17
$ws.Columns("X").ColumnWidth = 7.17
$ws.Columns("AA").ColumnWidth = 7.17
$ws.Columns("AB").ColumnWidth = 7.17
$ws.Columns("AC").ColumnWidth = 7.17
$ws.Columns("AD").ColumnWidth = 7.17
$ws.Columns("AH").ColumnWidth = 7.17